$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 42, shifting rows 42:62 down to 43:63.
$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with the new record's data.
$ws.Cells.Item(42, 1).Value = 4
$ws.Cells.Item(42, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(42, 3).Value = "Los Lagos"
$ws.Cells.Item(42, 4).Value = 44523
$ws.Cells.Item(42, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(42, 5).Value = 10
$ws.Cells.Item(42, 6).Value = 100112026
$ws.Cells.Item(42, 7).Value = "Haba"
$ws.Cells.Item(42, 8).Value = "Sin especificar"
$ws.Cells.Item(42, 9).Value = "Primera"
$ws.Cells.Item(42, 10).Value = 200
$ws.Cells.Item(42, 11).Value = 12000
$ws.Cells.Item(42, 12).Value = 12000
$ws.Cells.Item(42, 13).Value = 12000
$ws.Cells.Item(42, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(42, 15).Value = "Región del Maule"
$ws.Cells.Item(42, 16).Value = 480
$ws.Cells.Item(42, 17).Value = 25
$ws.Cells.Item(42, 18).Value = "Hortaliza"
